# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy the style from the existing header (G1) so H1 matches
# the bold/bordered/centered look of the other headers, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data rows: fill H2:H11 with 0.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
